$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to Text format so numeric-looking strings
# like "212.24" or "0.780" are stored verbatim instead of being
# coerced into floating point numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.265.24'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.603.61'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '212.24'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '0.485'
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '0.0614'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '18.19'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('D12').Value = '1.825.12'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '1.604.20'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').Value = '0.513'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').Value = '26.246.76'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '61.28'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '203.14'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').Value = '9.28'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('D23').Value = '6.01'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +11.77%  '
$ws.Range('D25').Value = '144.65'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -7.43%  '
$ws.Range('D28').Value = '15.17'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').Value = '6.55'
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('D30').Value = '0.0493'
$ws.Range('E30').Value = '  +3.77%  '
$ws.Range('D31').Value = '1.17'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').Value = '3.17'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = '2.92'
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('D34').Value = '2.42'
$ws.Range('E34').Value = '  +2.92%  '
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').Value = '1.148.24'
$ws.Range('E36').Value = '  +3.67%  '
$ws.Range('D37').Value = '0.0163'
$ws.Range('E37').Value = '  +6.65%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('D40').Value = '0.788'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').Value = '0.498'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('D42').Value = '0.780'
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').Value = '5.22'
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('D44').Value = '1.738.34'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').Value = '91.89'
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D46').Value = '1.52'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '54.17'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.407'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0957'
$ws.Range('E50').Value = '  -9.32%  '
$ws.Range('E51').Value = '  -0.03%  '

# Restore the default (unstyled) cell style so formatting matches
# the original workbook, which had no explicit style on these cells.
$priceRange.Style = "Normal"
